$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3131  # was 3129
$ws.Range("F4").Value = 1957  # was 1956
$ws.Range("F7").Value = 2622  # was 2620
$ws.Range("F13").Value = 140  # was 139
$ws.Range("F14").Value = 9865  # was 9856
$ws.Range("F16").Value = 25  # was 24
$ws.Range("F18").Value = 7796  # was 7788
$ws.Range("F19").Value = 12379  # was 12370
$ws.Range("F24").Value = 585  # was 583
$ws.Range("F25").Value = 2783  # was 2781
$ws.Range("F26").Value = 250  # was 249
$ws.Range("F27").Value = 226  # was 225
$ws.Range("F28").Value = 2823  # was 2818
$ws.Range("F29").Value = 1332  # was 1327
$ws.Range("F30").Value = 204  # was 60
$ws.Range("F32").Value = 70  # was 69
$ws.Range("F33").Value = 4596  # was 4593
$ws.Range("F34").Value = 1278  # was 1274
$ws.Range("F35").Value = 59  # was 58

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 645  # was 644

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 645  # was 644
$ws.Range("F4").Value = 3131  # was 3129
$ws.Range("F6").Value = 1957  # was 1956
$ws.Range("F9").Value = 2622  # was 2620
$ws.Range("F16").Value = 140  # was 139
$ws.Range("F17").Value = 9865  # was 9857
$ws.Range("F19").Value = 25  # was 24
$ws.Range("F21").Value = 7796  # was 7788
$ws.Range("F22").Value = 12379  # was 12370
$ws.Range("F27").Value = 585  # was 583
$ws.Range("F29").Value = 2783  # was 2781
$ws.Range("F32").Value = 250  # was 249
$ws.Range("F33").Value = 226  # was 225
$ws.Range("F34").Value = 206  # was 60
$ws.Range("F36").Value = 70  # was 69
$ws.Range("F37").Value = 4596  # was 4593

